# Add a "Countries" column (G) to the "Gem characteristics" sheet, fill it
# with "Unknown" plus a list of countries, format it like the neighbouring
# data column, sort the country list alphabetically, and leave the
# selection on G2 (matching the target edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell G1: text + formatting copied from A1 (the other header cells)
$ws.Range("G1").Value = "Countries"
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# G2 is fixed as "Unknown" (kept out of the alphabetical sort below)
$ws.Range("G2").Value = "Unknown"

# Countries, typed in this order (this becomes their insertion order in the
# shared-strings table); they get alphabetically sorted right after.
$countries = @(
    "Afghanistan",
    "Madagascar",
    "Pakistan",
    "Kenya",
    "USA",
    "Sri-lanka",
    "Nigeria",
    "Mozambique",
    "Brazil",
    "Australia",
    "Namibia",
    "Zambia",
    "Tanzania",
    "Birmania",
    "Thailand",
    "India",
    "Russia",
    "Colombia"
)

for ($i = 0; $i -lt $countries.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 7).Value = $countries[$i]
}

# Copy the data-cell formatting (same style as the rest of column F) onto
# the whole new column G range (header excluded).
$ws.Range("F2").Copy()
$ws.Range("G2:G20").PasteSpecial(-4122)

# Sort only the typed country list (G3:G20) alphabetically, leaving the
# "Unknown" entry in G2 untouched.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("G3:G20"))
$ws.Sort.SetRange($ws.Range("G3:G20"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Match the final selection recorded in the workbook.
$ws.Range("G2").Select()
